$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.859.61"
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = "'3.782.24"
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'600.52"
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').Value = "'164.78"
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').Value = "'0.158"
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').Value = "'0.450"
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('D11').Value = "'6.45"
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').Value = "'0.0000248"
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').Value = "'35.61"
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = "'4.426.38"
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = "'3.799.41"
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = "'67.921.52"
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = "'18.30"
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').Value = "'7.06"
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').Value = "'461.29"
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = "'9.69"
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('D22').Value = "'0.696"
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = "'0.0000148"
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').Value = "'82.69"
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = "'12.02"
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = "'9.98"
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = "'3.932.67"
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('D30').Value = "'7.40"
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = "'2.63"
$ws.Range('E31').Value = '  -5.33%  '
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').Value = "'29.21"
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = "'8.96"
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('D36').Value = "'0.0994"
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = "'3.24"
$ws.Range('E38').Value = '  -3.38%  '
$ws.Range('D39').Value = "'5.76"
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = "'47.37"
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').Value = "'43.07"
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('D46').Value = "'150.73"
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = "'8.33"
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +2.36%  '
$ws.Range('D49').Value = "'392.18"
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').Value = "'26.89"
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +6.75%  '
